# Generate Report for handback
# Refresh the handoff/handback datetimes for the "879e6e22-..." row (row 3)
# on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-01-18 02:51:40"
$wsZh.Range("G3").Value = "2016-01-18 02:52:22"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-01-18 02:51:50"
$wsDe.Range("G3").Value = "2016-01-18 02:52:39"
